# Apply the "JTS data mods" edit: rename survey sheets for parallel structure
# with the data-processing R script, and add two new metadata/ToDo sheets.

$wb = $excel.ActiveWorkbook

# --- Rename existing sheets ---------------------------------------------
# 1881            -> 1880Survey   (actually surveyed in 1881, tab kept as
#                                  "1880" for consistency with the R script)
# 1881notes       -> unchanged
# 1942            -> 1940Survey
# 1942notes       -> 1940notes
$wb.Worksheets.Item("1881").Name = "1880Survey"
$wb.Worksheets.Item("1942").Name = "1940Survey"
$wb.Worksheets.Item("1942notes").Name = "1940notes"

# --- Add the two new ToDo/metadata sheets at the end ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metadata1880 = $wb.Worksheets.Add($null, $lastSheet)
$metadata1880.Name = "1880Metadata"
$metadata1880.Range("A1").Value = "Actually surveyed in 1881; tab reads 1880 for consistency in the data processing R script"
$metadata1880.Range("A1").Font.Color = 0

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metadata1940 = $wb.Worksheets.Add($null, $lastSheet)
$metadata1940.Name = "1940Metadata"
$metadata1940.Range("A1").Value = "Kara please confirm this was actually sampled in 1940? Tab said 1942.."

# --- Make the 1940Survey tab the active one ------------------------------
$wb.Worksheets.Item("1940Survey").Activate()
